$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "RMB on object sprite" paragraph: merge the two split runs into one by
#    replacing the concatenated text with itself (Find/Replace coalesces the
#    run boundary that fell in the middle of "nothing").
# ---------------------------------------------------------------------------
$objSpriteOld = " " + [char]0x2013 + " Display any information available about the tile (or nothing if there is nothin" + "g unique about the tile)"
$rng = $d.Content
$rng.Find.Execute($objSpriteOld, $false, $false, $false, $false, $false, $true, 1, $false, $objSpriteOld, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Double-Click LMB on a tile/object/unit" paragraph.
#    a) remove the bold "/object/unit" run
#    b) rewrite the trailing description as four runs with the new wording
# ---------------------------------------------------------------------------
$dash = [char]0x2013

$rngTileP = $d.Content
$rngTileP.Find.Execute("/object/unit", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$oldDesc = " " + $dash + " Same effect as if RMB is clicked on it, as well as focusing the camera on the place that you double-clicked"
$rngDesc = $d.Content
$found = $rngDesc.Find.Execute($oldDesc, $false, $false, $false, $false, $false, $true, 1, $false, "", 1)
Write-Host "found desc to delete: " $found

Write-Host "tile paragraph text now:" $d.Paragraphs(18).Range.Text

# Append the four new plain-text runs right before the paragraph mark.
$endTile = $d.Paragraphs(18).Range
$endTile.Collapse(0) | Out-Null
$endTile.MoveEnd(1, -1) | Out-Null
$endTile.InsertAfter(" " + $dash + " Automatically move to the clicked tile")
$endTile.Collapse(0) | Out-Null
$endTile.InsertAfter(". This uses your move action")
$endTile.Collapse(0) | Out-Null
$endTile.InsertAfter(". ")
$endTile.Collapse(0) | Out-Null
$endTile.InsertAfter("The same also applies if you drag a custom path and then click quickly after release.")

Write-Host "tile paragraph final text:" $d.Paragraphs(18).Range.Text

# ---------------------------------------------------------------------------
# 3) Insert two brand-new paragraphs right after paragraph 18 ("... tile"):
#      - "Double-Click LMB on an interactable object - ..."
#      - "Double-Click LMB on an enemy - ..."
#    both indented with a 720-twip (36pt) first-line indent, matching the
#    other bulleted entries further down, and move the _GoBack bookmark into
#    the first of the two.
# ---------------------------------------------------------------------------

# --- paragraph: interactable object -----------------------------------
$tail = $d.Paragraphs(18).Range
$tail.Collapse(0) | Out-Null
$newParaA = $tail.InsertParagraphAfter()
$paraA = $d.Paragraphs(19).Range
$paraA.ParagraphFormat.FirstLineIndent = 36

$cur = $d.Paragraphs(19).Range
$cur.Collapse(0) | Out-Null
$cur.MoveEnd(1, -1) | Out-Null

$cur.InsertAfter("Double-Click LMB on an ")
$cur.Font.Bold = 1
$cur.Collapse(0) | Out-Null

$cur.InsertAfter("interactable")
$cur.Font.Bold = 1
$cur.Collapse(0) | Out-Null

$cur.InsertAfter(" object ")
$cur.Font.Bold = 1
$cur.Collapse(0) | Out-Null

$cur.InsertAfter($dash + " Automatically move to ")
$cur.Font.Bold = 0
$cur.Collapse(0) | Out-Null

$cur.InsertAfter("the closest adjacent tile and interact with the object")
$cur.Font.Bold = 0
$cur.Collapse(0) | Out-Null

$cur.InsertAfter(". ")
$cur.Font.Bold = 0
$cur.Collapse(0) | Out-Null

$cur.InsertAfter("This uses your move (if you need to move there it interact with it) ")
$cur.Font.Bold = 0
$cur.Collapse(0) | Out-Null

# --- move the _GoBack bookmark here (zero-width, between the two runs) --
$bmSpot = $cur.Duplicate
$bmSpot.InsertAfter("#")
$bmSpot.Font.Bold = 0
$d.Bookmarks.Add("_GoBack", $bmSpot)
$bmSpot.Delete()
$cur.Collapse(0) | Out-Null

$cur.InsertAfter("and a minor action. ")
$cur.Font.Bold = 0
$cur.Collapse(0) | Out-Null

$cur.InsertAfter("The same also applies if you drag a custom path and then click quickly after release.")
$cur.Font.Bold = 0
$cur.Collapse(0) | Out-Null

Write-Host "interactable paragraph text:" $d.Paragraphs(19).Range.Text

